{"js": "// Locate the relevant paragraphs by their placeholder text, then:\n//  1) Merge the PERSON_77 / PERSON_78 paragraphs into a single paragraph\n//     reading \"[[PERSON_77]] - [[PERSON_78]], [[PERSON_77]]\" and drop the\n//     now-redundant PERSON_78 paragraph entirely.\n//  2) Append a brand-new list paragraph for PERSON_102 right after the\n//     PERSON_101 paragraph (before the trailing empty paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet p77 = null;\nlet p78 = null;\nlet p101 = null;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (p77 === null && t === \"[[PERSON_77]] \u2013 [[PERSON_77]], [[PERSON_77]]\") {\n    p77 = items[i];\n  } else if (p78 === null && t === \"[[PERSON_78]] \u2013 [[PERSON_78]], [[PERSON_78]]\") {\n    p78 = items[i];\n  } else if (t === \"[[PERSON_101]] \u2013 [[PERSON_101]], [[PERSON_101]]\") {\n    p101 = items[i];\n  }\n}\n\nif (!p77 || !p78 || !p101) {\n  throw new Error(\n    \"Could not locate expected paragraphs (p77=\" + !!p77 +\n    \", p78=\" + !!p78 + \", p101=\" + !!p101 + \")\"\n  );\n}\n\n// 1a) Rewrite PERSON_77's paragraph text to fold in PERSON_78's name.\np77.insertText(\"[[PERSON_77]] \u2013 [[PERSON_78]], [[PERSON_77]]\", \"Replace\");\n\n// 1b) Remove the (now duplicate) PERSON_78 paragraph entirely.\np78.delete();\n\n// 2) Insert a fresh PERSON_102 list paragraph right after PERSON_101.\np101.insertParagraph(\"[[PERSON_102]] \u2013 [[PERSON_102]], [[PERSON_102]]\", \"After\");\n\nawait context.sync();\n", "ps1": "# Locate the relevant paragraphs by their placeholder text, then:\n#  1) Merge the PERSON_77 / PERSON_78 paragraphs into a single paragraph\n#     reading \"[[PERSON_77]] - [[PERSON_78]], [[PERSON_77]]\" and drop the\n#     now-redundant PERSON_78 paragraph entirely.\n#  2) Append a brand-new list paragraph for PERSON_102 right after the\n#     PERSON_101 paragraph (before the trailing empty paragraph).\n\n$d = $word.ActiveDocument\n\n$target77 = \"[[PERSON_77]] \" + [char]8211 + \" [[PERSON_77]], [[PERSON_77]]\"\n$target78 = \"[[PERSON_78]] \" + [char]8211 + \" [[PERSON_78]], [[PERSON_78]]\"\n$target101 = \"[[PERSON_101]] \" + [char]8211 + \" [[PERSON_101]], [[PERSON_101]]\"\n\n$p77 = $null\n$p78 = $null\n$p101 = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $target77) {\n        $p77 = $p\n    } elseif ($t -eq $target78) {\n        $p78 = $p\n    } elseif ($t -eq $target101) {\n        $p101 = $p\n    }\n}\n\nif ($p77 -eq $null -or $p78 -eq $null -or $p101 -eq $null) {\n    Write-Output \"ERROR: could not locate expected paragraphs\"\n} else {\n    # 1a) Rewrite PERSON_77's paragraph text to fold in PERSON_78's name.\n    $p77.Range.Text = \"[[PERSON_77]] \" + [char]8211 + \" [[PERSON_78]], [[PERSON_77]]\"\n\n    # 1b) Remove the (now duplicate) PERSON_78 paragraph entirely. This\n    #     shifts every later paragraph's index down by one, so any\n    #     previously-cached paragraph reference/index past this point is\n    #     stale and must be re-resolved (done below for PERSON_101).\n    $p78.Range.Delete()\n\n    # Re-locate the PERSON_101 paragraph now that indices have shifted.\n    $p101 = $null\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $target101) {\n            $p101 = $p\n        }\n    }\n\n    # 2) Insert a fresh PERSON_102 list paragraph right after PERSON_101.\n    $p101.Range.InsertParagraphAfter()\n    $newIndex = $p101.Index + 1\n    $newp = $d.Paragraphs.Item($newIndex)\n    $newp.Range.Text = \"[[PERSON_102]] \" + [char]8211 + \" [[PERSON_102]], [[PERSON_102]]\"\n}\n"}
